$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.366566
$ws.Range("H2").Value = 1.099698
$ws.Range("I2").Value = 0.5689653834353526
$ws.Range("J2").Value = 0.5689653834353527
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.3777296666666667
$ws.Range("N2").Value = 1.133189
$ws.Range("O2").Value = 0.2121489992374768
$ws.Range("P2").Value = 0.2121489992374768
$ws.Range("Q2").Value = 0.1384628529913333
$ws.Range("R2").Value = 1.246165676922
$ws.Range("S2").Value = 0.1207054366965774
$ws.Range("T2").Value = 0.1207054366965774

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.366566
$ws.Range("H3").Value = 1.099698
$ws.Range("I3").Value = 0.5689653834353526
$ws.Range("J3").Value = 0.5689653834353527
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.092289666666667
$ws.Range("N3").Value = 3.276869
$ws.Range("O3").Value = 0.6134761976883921
$ws.Range("P3").Value = 0.6134761976883921
$ws.Range("Q3").Value = 0.4003962539513333
$ws.Range("R3").Value = 3.603566285562
$ws.Range("S3").Value = 0.3490467200462382
$ws.Range("T3").Value = 0.3490467200462383

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.366566
$ws.Range("H4").Value = 1.099698
$ws.Range("I4").Value = 0.5689653834353526
$ws.Range("J4").Value = 0.5689653834353527
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.310473
$ws.Range("N4").Value = 0.931419
$ws.Range("O4").Value = 0.174374803074131
$ws.Range("P4").Value = 0.174374803074131
$ws.Range("Q4").Value = 0.113808845718
$ws.Range("R4").Value = 1.024279611462
$ws.Range("S4").Value = 0.09921322669253706
$ws.Range("T4").Value = 0.09921322669253706

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.2777016666666667
$ws.Range("H5").Value = 0.833105
$ws.Range("I5").Value = 0.4310346165646473
$ws.Range("J5").Value = 0.4310346165646473
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.3777296666666667
$ws.Range("N5").Value = 1.133189
$ws.Range("O5").Value = 0.2121489992374768
$ws.Range("P5").Value = 0.2121489992374768
$ws.Range("Q5").Value = 0.1048961579827778
$ws.Range("R5").Value = 0.944065421845
$ws.Range("S5").Value = 0.09144356254089948
$ws.Range("T5").Value = 0.09144356254089947

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.2777016666666667
$ws.Range("H6").Value = 0.833105
$ws.Range("I6").Value = 0.4310346165646473
$ws.Range("J6").Value = 0.4310346165646473
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.092289666666667
$ws.Range("N6").Value = 3.276869
$ws.Range("O6").Value = 0.6134761976883921
$ws.Range("P6").Value = 0.6134761976883921
$ws.Range("Q6").Value = 0.3033306609161111
$ws.Range("R6").Value = 2.729975948245
$ws.Range("S6").Value = 0.2644294776421539
$ws.Range("T6").Value = 0.2644294776421539

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.2777016666666667
$ws.Range("H7").Value = 0.833105
$ws.Range("I7").Value = 0.4310346165646473
$ws.Range("J7").Value = 0.4310346165646473
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.310473
$ws.Range("N7").Value = 0.931419
$ws.Range("O7").Value = 0.174374803074131
$ws.Range("P7").Value = 0.174374803074131
$ws.Range("Q7").Value = 0.08621886955500001
$ws.Range("R7").Value = 0.775969825995
$ws.Range("S7").Value = 0.07516157638159394
$ws.Range("T7").Value = 0.07516157638159393

